# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# --- Bad Drivers table ---
# Row 3: Intel(R) Wi-Fi 6E AX210 160MHz - 23.60.1.2
$ws.Range("C3").Value = 619
$ws.Range("D3").Value = 86.8

# Row 5: Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.0.4
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 254
$ws.Range("D5").Value = 98.9

# Row 6: Totals
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 1178

# --- Good Drivers table ---
# Row 16: Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Range("B16").Value = 449371

# Row 17: MediaTek Wi-Fi 6E MT7922 (RZ616) 160MHz PCIe Adapter - 3.4.0.1088
$ws.Range("B17").Value = 86281

# Row 19: Intel(R) Wi-Fi 6E AX210 160MHz - 23.20.1.1
$ws.Range("B19").Value = 14968

# Row 24: Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9
$ws.Range("B24").Value = 77999
